# Auto-generated Excel COM-interop edit script
# Applies scheduled-runner market-data updates to the Cactuar_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H111").Value = 1733.4286
$ws.Range("I111").Value = 1567.125
$ws.Range("K111").Value = 4701.375
$ws.Range("M111").Value = -1634.375
$ws.Range("H112").Value = 3248.375
$ws.Range("I112").Value = 1395
$ws.Range("K112").Value = 4185
$ws.Range("M112").Value = -3077
$ws.Range("H132").Value = 7851.756
$ws.Range("I132").Value = 2205.4146
$ws.Range("J132").Value = 12996.2
$ws.Range("K132").Value = 6616.2438
$ws.Range("L132").Value = 38988.60000000001
$ws.Range("M132").Value = -4086.2438
$ws.Range("N132").Value = -44048.60000000001
$ws.Range("H137").Value = 6947768.5
$ws.Range("I137").Value = 1077.2693
$ws.Range("K137").Value = 3231.8079
$ws.Range("M137").Value = -681.8078999999998
$ws.Range("H138").Value = 3871
$ws.Range("I138").Value = 2308.125
$ws.Range("J138").Value = 4175.951
$ws.Range("K138").Value = 6924.375
$ws.Range("L138").Value = 12527.853
$ws.Range("M138").Value = -1784.375
$ws.Range("N138").Value = -22807.853
$ws.Range("H139").Value = 123666.336
$ws.Range("J139").Value = 123666.336
$ws.Range("L139").Value = 123666.336
$ws.Range("N139").Value = -133946.336
$ws.Range("H140").Value = 57163.8
$ws.Range("J140").Value = 55658.777
$ws.Range("L140").Value = 55658.777
$ws.Range("N140").Value = -66018.777
$ws.Range("H141").Value = 5866.607
$ws.Range("I141").Value = 5164.3076
$ws.Range("K141").Value = 15492.9228
$ws.Range("M141").Value = -10312.9228

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1247770.9
$ws.Range("I2").Value = 1939711.1
$ws.Range("J2").Value = 2278.6
$ws.Range("K2").Value = 1939711.1
$ws.Range("L2").Value = 2278.6
$ws.Range("M2").Value = -1939598.1
$ws.Range("N2").Value = -2504.6
$ws.Range("H32").Value = 4305.557
$ws.Range("I32").Value = 2308.8728
$ws.Range("K32").Value = 2308.8728
$ws.Range("M32").Value = -2021.8728
$ws.Range("H110").Value = 819195.1
$ws.Range("I110").Value = 1075378
$ws.Range("K110").Value = 1075378
$ws.Range("M110").Value = -1073333
$ws.Range("H116").Value = 1247770.9
$ws.Range("I116").Value = 1939711.1
$ws.Range("J116").Value = 2278.6
$ws.Range("K116").Value = 1939711.1
$ws.Range("L116").Value = 2278.6
$ws.Range("M116").Value = -1937417.1
$ws.Range("N116").Value = -6866.6
$ws.Range("H122").Value = 3851.2354
$ws.Range("I122").Value = 2486.4546
$ws.Range("J122").Value = 6353.3335
$ws.Range("K122").Value = 7459.3638
$ws.Range("L122").Value = 19060.0005
$ws.Range("M122").Value = -5009.3638
$ws.Range("N122").Value = -23960.0005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1247770.9
$ws.Range("I3").Value = 1939711.1
$ws.Range("J3").Value = 2278.6
$ws.Range("K3").Value = 1939711.1
$ws.Range("L3").Value = 2278.6
$ws.Range("M3").Value = -1939597.1
$ws.Range("N3").Value = -2506.6
$ws.Range("H13").Value = 70999
$ws.Range("I13").Value = 70999
$ws.Range("K13").Value = 70999
$ws.Range("M13").Value = -70831
$ws.Range("H107").Value = 1376.9032
$ws.Range("I107").Value = 1297.2
$ws.Range("J107").Value = 1709
$ws.Range("K107").Value = 1297.2
$ws.Range("L107").Value = 1709
$ws.Range("M107").Value = 622.8
$ws.Range("N107").Value = -5549

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3132.0952
$ws.Range("I31").Value = 999
$ws.Range("K31").Value = 999
$ws.Range("M31").Value = -704
$ws.Range("H34").Value = 3132.0952
$ws.Range("I34").Value = 999
$ws.Range("K34").Value = 999
$ws.Range("M34").Value = -797
$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 35000
$ws.Range("L97").Value = 35000
$ws.Range("N97").Value = -36982
$ws.Range("H141").Value = 92179
$ws.Range("I141").Value = 39494
$ws.Range("J141").Value = 100959.836
$ws.Range("K141").Value = 39494
$ws.Range("L141").Value = 100959.836
$ws.Range("M141").Value = -34314
$ws.Range("N141").Value = -111319.836

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 12823488
$ws.Range("I131").Value = 55556170
$ws.Range("J131").Value = 9262432
$ws.Range("K131").Value = 166668510
$ws.Range("L131").Value = 27787296
$ws.Range("M131").Value = -166663470
$ws.Range("N131").Value = -27797376

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 10855407
$ws.Range("I102").Value = 13080918
$ws.Range("J102").Value = 6036.625
$ws.Range("K102").Value = 13080918
$ws.Range("L102").Value = 6036.625
$ws.Range("M102").Value = -13079296
$ws.Range("N102").Value = -9280.625
$ws.Range("H113").Value = 3220
$ws.Range("I113").Value = 2700
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -530
$ws.Range("N113").Value = -8340
$ws.Range("H126").Value = 3741.8928
$ws.Range("I126").Value = 3024.4375
$ws.Range("K126").Value = 9073.3125
$ws.Range("M126").Value = -6603.3125
$ws.Range("H135").Value = 100780
$ws.Range("J135").Value = 100780
$ws.Range("L135").Value = 100780
$ws.Range("N135").Value = -110920
$ws.Range("H140").Value = 78699
$ws.Range("J140").Value = 78699
$ws.Range("L140").Value = 78699
$ws.Range("N140").Value = -89059

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H68").Value = 912773.7
$ws.Range("I68").Value = 1422742.5
$ws.Range("J68").Value = 6162.5557
$ws.Range("K68").Value = 1422742.5
$ws.Range("L68").Value = 6162.5557
$ws.Range("M68").Value = -1421993.5
$ws.Range("N68").Value = -7660.5557
$ws.Range("H71").Value = 912773.7
$ws.Range("I71").Value = 1422742.5
$ws.Range("J71").Value = 6162.5557
$ws.Range("K71").Value = 7113712.5
$ws.Range("L71").Value = 30812.7785
$ws.Range("M71").Value = -7109968.5
$ws.Range("N71").Value = -38300.7785
$ws.Range("H93").Value = 3736.5
$ws.Range("I93").Value = 979.5
$ws.Range("K93").Value = 979.5
$ws.Range("M93").Value = 268.5
$ws.Range("H100").Value = 3758.8
$ws.Range("I100").Value = 3897.5
$ws.Range("K100").Value = 3897.5
$ws.Range("M100").Value = -3356.5
$ws.Range("H132").Value = 4057.5352
$ws.Range("I132").Value = 3109.6875
$ws.Range("K132").Value = 9329.0625
$ws.Range("M132").Value = -6799.0625
$ws.Range("H139").Value = 88358
$ws.Range("J139").Value = 88358
$ws.Range("L139").Value = 88358
$ws.Range("N139").Value = -98638
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 120000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 120000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 120000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -121872
$ws.Range("H78").Value = 120000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 120000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 360000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -369360
$ws.Range("H107").Value = 2788.375
$ws.Range("I107").Value = 3356.5
$ws.Range("J107").Value = 1084
$ws.Range("K107").Value = 10069.5
$ws.Range("L107").Value = 3252
$ws.Range("M107").Value = -8149.5
$ws.Range("N107").Value = -7092
$ws.Range("H118").Value = 99245.5
$ws.Range("J118").Value = 99245.5
$ws.Range("L118").Value = 99245.5
$ws.Range("N118").Value = -102559.5
$ws.Range("H138").Value = 78000
$ws.Range("J138").Value = 78000
$ws.Range("L138").Value = 78000
$ws.Range("N138").Value = -88280

